# Updates cryptos list figures (price/volume) to match the latest scrape,
# including a few coin-rank swaps (rows 34/35, 40/42, 50/51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.862.25'
$ws.Range('E2').Value = '  -5.61%  '

$ws.Range('D3').Value = '2.993.25'
$ws.Range('E3').Value = '  -6.49%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '548.70'
$ws.Range('E5').Value = '  -4.96%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.69'
$ws.Range('E6').Value = '  -8.53%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.565'
$ws.Range('E8').Value = '  -5.47%  '

$ws.Range('D9').Value = '2.995.79'
$ws.Range('E9').Value = '  -6.19%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.111'
$ws.Range('E10').Value = '  -6.96%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.22'
$ws.Range('E11').Value = '  -7.44%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.363'
$ws.Range('E12').Value = '  -6.74%  '

$ws.Range('D13').Value = '3.512.32'
$ws.Range('E13').Value = '  -6.44%  '

$ws.Range('E14').Value = '  -3.69%  '

$ws.Range('D15').Value = '61.915.92'
$ws.Range('E15').Value = '  -5.38%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.63'
$ws.Range('E16').Value = '  -7.97%  '

$ws.Range('D17').Value = '2.993.50'
$ws.Range('E17').Value = '  -6.08%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000147'
$ws.Range('E18').Value = '  -6.95%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '389.73'
$ws.Range('E19').Value = '  -6.02%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.10'
$ws.Range('E20').Value = '  -4.57%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.92'
$ws.Range('E21').Value = '  -7.47%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.62'
$ws.Range('E22').Value = '  -7.81%  '

$ws.Range('E23').Value = '  +0.08%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.92'
$ws.Range('E24').Value = '  -6.45%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.468'
$ws.Range('E25').Value = '  -4.56%  '

$ws.Range('E26').Value = '  -8.09%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.02%  '

$ws.Range('D28').Value = '0.0₃0936'
$ws.Range('E28').Value = '  -11.68%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.42'
$ws.Range('E29').Value = '  -6.11%  '

$ws.Range('E30').Value = '  -0.03%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.72'
$ws.Range('E31').Value = '  -6.89%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.43'
$ws.Range('E32').Value = '  -5.54%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '159.41'
$ws.Range('E33').Value = '  +1.83%  '

$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.00'
$ws.Range('E34').Value = '  -6.37%  '

$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.62'
$ws.Range('E35').Value = '  -8.13%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.07'
$ws.Range('E36').Value = '  -6.80%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.28'
$ws.Range('E37').Value = '  -6.79%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.56'
$ws.Range('E38').Value = '  -10.22%  '

$ws.Range('D39').Value = '2.433.34'
$ws.Range('E39').Value = '  -11.17%  '

$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.20'
$ws.Range('E40').Value = '  -4.66%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.88'
$ws.Range('E41').Value = '  -7.05%  '

$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.33'
$ws.Range('E42').Value = '  -7.92%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.660'
$ws.Range('E43').Value = '  -7.31%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0594'
$ws.Range('E44').Value = '  -6.81%  '

$ws.Range('E45').Value = '  -0.13%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0245'
$ws.Range('E46').Value = '  -7.35%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.92'
$ws.Range('E47').Value = '  -11.62%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0951'
$ws.Range('E48').Value = '  -4.02%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.66'
$ws.Range('E49').Value = '  -9.35%  '

$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.47'
$ws.Range('E50').Value = '  +0.18%  '

$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '264.54'
$ws.Range('E51').Value = '  -10.91%  '
